$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

# Row 2
Set-TextValue $ws.Range("D2") "36.949.14"
$ws.Range("E2").Value = "  -0.99%  "

# Row 3
Set-TextValue $ws.Range("D3") "2.000.64"
$ws.Range("E3").Value = "  -2.39%  "

# Row 4
Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.17%  "

# Row 5
$ws.Range("E5").Value = "  -1.72%  "

# Row 6
$ws.Range("E6").Value = "  -1.62%  "

# Row 7
$ws.Range("E7").Value = "  +0.00%  "

# Row 8
Set-TextValue $ws.Range("D8") "54.50"
$ws.Range("E8").Value = "  -4.72%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.373"
$ws.Range("E9").Value = "  -3.45%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.0777"
$ws.Range("E10").Value = "  -3.67%  "

# Row 11
$ws.Range("E11").Value = "  -5.66%  "

# Row 12
Set-TextValue $ws.Range("D12") "2.293.51"
$ws.Range("E12").Value = "  -2.47%  "

# Row 13
Set-TextValue $ws.Range("D13") "13.86"
$ws.Range("E13").Value = "  -5.71%  "

# Row 14
Set-TextValue $ws.Range("D14") "19.68"
$ws.Range("E14").Value = "  -5.41%  "

# Row 15
Set-TextValue $ws.Range("D15") "5.19"
$ws.Range("E15").Value = "  -2.33%  "

# Row 16
$ws.Range("E16").Value = "  -3.75%  "

# Row 17
Set-TextValue $ws.Range("D17") "2.007.91"
$ws.Range("E17").Value = "  -1.91%  "

# Row 18
Set-TextValue $ws.Range("D18") "36.871.41"
$ws.Range("E18").Value = "  -0.93%  "

# Row 19
Set-TextValue $ws.Range("D19") "6.27"
$ws.Range("E19").Value = "  +2.97%  "

# Row 20
Set-TextValue $ws.Range("D20") "68.16"
$ws.Range("E20").Value = "  -2.39%  "

# Row 21
Set-TextValue $ws.Range("D21") "0.0₃0809"
$ws.Range("E21").Value = "  -3.14%  "

# Row 22
Set-TextValue $ws.Range("D22") "221.24"
$ws.Range("E22").Value = "  -2.30%  "

# Row 23
$ws.Range("E23").Value = "  +0.11%  "

# Row 24
$ws.Range("E24").Value = "  +0.98%  "

# Row 25
$ws.Range("E25").Value = "  -6.10%  "

# Row 26
Set-TextValue $ws.Range("D26") "164.02"
$ws.Range("E26").Value = "  -2.54%  "

# Row 27
Set-TextValue $ws.Range("D27") "8.94"
$ws.Range("E27").Value = "  -6.67%  "

# Row 28
$ws.Range("E28").Value = "  -4.19%  "

# Row 29
Set-TextValue $ws.Range("D29") "18.47"
$ws.Range("E29").Value = "  -2.61%  "

# Row 30
$ws.Range("E30").Value = "  -7.26%  "

# Row 31
$ws.Range("E31").Value = "  -1.94%  "

# Row 32
$ws.Range("E32").Value = "  -2.51%  "

# Row 33
$ws.Range("E33").Value = "  -3.06%  "

# Row 34
Set-TextValue $ws.Range("D34") "4.44"
$ws.Range("E34").Value = "  -3.96%  "

# Row 35
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws.Range("D35") "2.30"
$ws.Range("E35").Value = "  -4.92%  "

# Row 36
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D36") "1.86"
$ws.Range("E36").Value = "  +1.68%  "

# Row 37
Set-TextValue $ws.Range("D37") "0.999"
$ws.Range("E37").Value = "  -0.19%  "

# Row 38
Set-TextValue $ws.Range("D38") "3.10"
$ws.Range("E38").Value = "  -4.57%  "

# Row 39
$ws.Range("E39").Value = "  -1.21%  "

# Row 40
Set-TextValue $ws.Range("D40") "1.450.76"
$ws.Range("E40").Value = "  -2.86%  "

# Row 41
Set-TextValue $ws.Range("D41") "0.0211"
$ws.Range("E41").Value = "  -5.18%  "

# Row 42
Set-TextValue $ws.Range("D42") "94.35"
$ws.Range("E42").Value = "  -2.38%  "

# Row 43
Set-TextValue $ws.Range("D43") "2.78"
$ws.Range("E43").Value = "  -4.22%  "

# Row 44
Set-TextValue $ws.Range("D44") "0.0903"
$ws.Range("E44").Value = "  -4.08%  "

# Row 45
$ws.Range("E45").Value = "  -4.55%  "

# Row 46
Set-TextValue $ws.Range("D46") "15.70"
$ws.Range("E46").Value = "  -8.62%  "

# Row 47
Set-TextValue $ws.Range("D47") "7.08"
$ws.Range("E47").Value = "  -1.13%  "

# Row 48
$ws.Range("E48").Value = "  -3.12%  "

# Row 49
Set-TextValue $ws.Range("D49") "2.88"
$ws.Range("E49").Value = "  -1.29%  "

# Row 50
Set-TextValue $ws.Range("D50") "2.182.88"
$ws.Range("E50").Value = "  -2.57%  "

# Row 51
Set-TextValue $ws.Range("D51") "3.54"
$ws.Range("E51").Value = "  -8.50%  "
